$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Rename "Requested quantity" headers on existing sheets
$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "PO Forecast"

# Copy cell formatting (font/border/alignment) from the "Weekly Quantity" header row
$ws1.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

# Copy the date-formatted style used in column A of the data rows
$ws1.Range("A2").Copy()
$newSheet.Range("A2:A55").PasteSpecial(-4122)

# Header row values
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# Data values
$colA = @(44955.99999999999,45018.99999999999,45032.99999999999,45067.99999999999,45074.99999999999,45081.99999999999,45088.99999999999,45095.99999999999,45102.99999999999,45109.99999999999,45123.99999999999,45130.99999999999,45137.99999999999,45165.99999999999,45235.99999999999,45249.99999999999,45256.99999999999,45263.99999999999,45298.99999999999,45305.99999999999,45312.99999999999,45333.99999999999,45340.99999999999,45347.99999999999,45354.99999999999,45361.99999999999,45368.99999999999,45382.99999999999,45389.99999999999,45396.99999999999,45403.99999999999,45410.99999999999,45417.99999999999,45424.99999999999,45445.99999999999,45494.99999999999,45501.99999999999,45515.99999999999,45522.99999999999,45529.99999999999,45550.99999999999,45557.99999999999,45564.99999999999,45571.99999999999,45585.99999999999,45599.99999999999,45606.99999999999,45613.99999999999,45620.99999999999,45627.99999999999,45634.99999999999,45641.99999999999,45648.99999999999,45655.99999999999)
$colB = @(70,75,76,79,79,80,80,81,81,82,83,83,84,86,91,92,93,94,96,97,97,99,99,100,100,101,101,103,103,104,104,105,105,106,107,111,112,113,113,114,115,116,116,117,118,119,119,120,120,121,122,122,123,123)
$colC = @(-77.95452413601268,-71.81240080393177,-79.62264268756921,-64.76621331287419,-68.20902786620114,-68.99847980947841,-62.54004041737574,-70.2131731027112,-72.0723141571908,-60.84784468483574,-64.00832359083297,-66.50542555795823,-67.99539233789451,-67.74139333209362,-71.4688908660289,-67.02723668655298,-50.79309597156892,-62.09015111359153,-44.86790429710968,-60.16562038156194,-51.74329431879529,-54.96463594719937,-46.73834418831726,-42.40924792353879,-44.83863307304108,-32.60759353390347,-51.84176878806007,-53.58583620129977,-41.84760568942556,-47.81971314739381,-42.68186167937606,-51.31586628689925,-50.49718844168113,-46.4150156458661,-48.85617610935928,-39.49804816297192,-34.73991124042663,-33.39138366977189,-40.19005178711286,-40.1572649941246,-46.12782990563502,-28.42527974557532,-34.11957412649119,-32.2524822471669,-30.31125390059152,-27.21459211204624,-37.82038541350573,-29.94755533064287,-27.48149862395706,-29.76788140239055,-25.03617790005428,-28.65297296299482,-33.94034751937073,-15.83410273026471)
$colD = @(224.3433738376704,232.772098618468,226.6877236290735,239.9513481728172,230.8555285464355,217.3302393342116,226.1111534327484,239.0787421245435,232.2104766275146,234.9951928929652,220.103876714955,228.3559507647531,239.4297698536421,227.8054910772229,232.1094529363405,242.2775577909624,249.1480707946022,240.2832758873422,255.0211998346596,249.8512209806746,245.0556616141128,245.9838640972291,248.9148451613787,246.5550787408174,250.8404242340377,255.1575693579658,251.3890920205423,244.3160510276288,255.9323084938668,257.7090433478691,265.3322861860874,250.8251494868412,251.8415433020244,262.4241274131418,254.1620077155121,274.7295843483599,259.4550211292342,273.1265002263449,264.4746623138895,266.0856660386406,261.1077447886579,271.4270328818843,277.1785355288732,255.4149612226694,269.3181647113717,270.245339936985,264.8666697744933,279.1921629520487,281.8372412451992,277.1762616495807,284.1203763640547,269.4652550918231,268.4256200574456,281.9192284755436)

for ($i = 0; $i -lt $colA.Length; $i++) {
    $r = $i + 2
    $newSheet.Cells.Item($r, 1).Value = $colA[$i]
    $newSheet.Cells.Item($r, 2).Value = $colB[$i]
    $newSheet.Cells.Item($r, 3).Value = $colC[$i]
    $newSheet.Cells.Item($r, 4).Value = $colD[$i]
}

Write-Host "PO Forecast sheet populated"
